$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): shared-string/column ordering was restructured ---
$ws.Cells.Item(1, 1).Value = "L"
$ws.Cells.Item(1, 2).Value = "Acetic acid43"
$ws.Cells.Item(1, 3).Value = "Water"
$ws.Cells.Item(1, 4).Value = "H2"
$ws.Cells.Item(1, 5).Value = "Ketene"
$ws.Cells.Item(1, 6).Value = "CO"
$ws.Cells.Item(1, 7).Value = "Acetaldehyde"
$ws.Cells.Item(1, 8).Value = "Acetic acid60"
$ws.Cells.Item(1, 9).Value = "CO2"

# --- Write the updated/restructured data block (rows 2-9 overwritten in place) ---
$ws.Cells.Item(2, 1).Value = 0.0
$ws.Cells.Item(2, 2).Value = -1.0
$ws.Cells.Item(2, 3).Value = -1.0
$ws.Cells.Item(2, 4).Value = 0.000001345079964297693
$ws.Cells.Item(2, 5).Value = -1.0
$ws.Cells.Item(2, 6).Value = 0.0000007854158387545218
$ws.Cells.Item(2, 7).Value = -1.0
$ws.Cells.Item(2, 8).Value = -1.0
$ws.Cells.Item(2, 9).Value = 0.000000151051300016021

$ws.Cells.Item(3, 1).Value = 0.0
$ws.Cells.Item(3, 2).Value = -1.0
$ws.Cells.Item(3, 3).Value = -1.0
$ws.Cells.Item(3, 4).Value = 0.000001424790176625913
$ws.Cells.Item(3, 5).Value = -1.0
$ws.Cells.Item(3, 6).Value = 0.0000009668408433756081
$ws.Cells.Item(3, 7).Value = -1.0
$ws.Cells.Item(3, 8).Value = -1.0
$ws.Cells.Item(3, 9).Value = 0.0000001545272978201269

$ws.Cells.Item(4, 1).Value = 0.0004
$ws.Cells.Item(4, 2).Value = -1.0
$ws.Cells.Item(4, 3).Value = -1.0
$ws.Cells.Item(4, 4).Value = 0.000001070725567250404
$ws.Cells.Item(4, 5).Value = -1.0
$ws.Cells.Item(4, 6).Value = 0.0000003609250316013533
$ws.Cells.Item(4, 7).Value = -1.0
$ws.Cells.Item(4, 8).Value = -1.0
$ws.Cells.Item(4, 9).Value = 0.00000006885662223673767

$ws.Cells.Item(5, 1).Value = 0.0008
$ws.Cells.Item(5, 2).Value = -1.0
$ws.Cells.Item(5, 3).Value = -1.0
$ws.Cells.Item(5, 4).Value = 0.000000971629297488102
$ws.Cells.Item(5, 5).Value = -1.0
$ws.Cells.Item(5, 6).Value = 0.0000005344902859832635
$ws.Cells.Item(5, 7).Value = -1.0
$ws.Cells.Item(5, 8).Value = -1.0
$ws.Cells.Item(5, 9).Value = 0.00000008347650509387418

$ws.Cells.Item(6, 1).Value = 0.0012
$ws.Cells.Item(6, 2).Value = -1.0
$ws.Cells.Item(6, 3).Value = -1.0
$ws.Cells.Item(6, 4).Value = 0.0000008498638856125455
$ws.Cells.Item(6, 5).Value = -1.0
$ws.Cells.Item(6, 6).Value = 0.0000005684576476198015
$ws.Cells.Item(6, 7).Value = -1.0
$ws.Cells.Item(6, 8).Value = -1.0
$ws.Cells.Item(6, 9).Value = 0.00000008346556882406874

$ws.Cells.Item(7, 1).Value = 0.0015
$ws.Cells.Item(7, 2).Value = -1.0
$ws.Cells.Item(7, 3).Value = -1.0
$ws.Cells.Item(7, 4).Value = 0.000001291442370793711
$ws.Cells.Item(7, 5).Value = -1.0
$ws.Cells.Item(7, 6).Value = 0.0000007130362430638149
$ws.Cells.Item(7, 7).Value = -1.0
$ws.Cells.Item(7, 8).Value = -1.0
$ws.Cells.Item(7, 9).Value = 0.0000001174080074705969

$ws.Cells.Item(8, 1).Value = 0.003
$ws.Cells.Item(8, 2).Value = -1.0
$ws.Cells.Item(8, 3).Value = -1.0
$ws.Cells.Item(8, 4).Value = 0.000002102430672936632
$ws.Cells.Item(8, 5).Value = -1.0
$ws.Cells.Item(8, 6).Value = 0.000001393567697711225
$ws.Cells.Item(8, 7).Value = -1.0
$ws.Cells.Item(8, 8).Value = -1.0
$ws.Cells.Item(8, 9).Value = 0.00000015506199015438

$ws.Cells.Item(9, 1).Value = 0.004
$ws.Cells.Item(9, 2).Value = -1.0
$ws.Cells.Item(9, 3).Value = -1.0
$ws.Cells.Item(9, 4).Value = 0.0000009732524784643293
$ws.Cells.Item(9, 5).Value = -1.0
$ws.Cells.Item(9, 6).Value = 0.0000007279948254386134
$ws.Cells.Item(9, 7).Value = -1.0
$ws.Cells.Item(9, 8).Value = -1.0
$ws.Cells.Item(9, 9).Value = 0.0000001059777857182084

# --- Add new row 10 (0.005 data point), shifted down from the old row 9 ---
$ws.Cells.Item(10, 1).Value = 0.005
$ws.Cells.Item(10, 2).Value = -1.0
$ws.Cells.Item(10, 3).Value = -1.0
$ws.Cells.Item(10, 4).Value = 0.000002262929135699318
$ws.Cells.Item(10, 5).Value = -1.0
$ws.Cells.Item(10, 6).Value = 0.000001963475081997287
$ws.Cells.Item(10, 7).Value = -1.0
$ws.Cells.Item(10, 8).Value = -1.0
$ws.Cells.Item(10, 9).Value = 0.0000002289070786371638

# Copy column-A formatting (border/bold/alignment) down into the new row 10 cell
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = 0
